$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 7 data, mirroring the style/format of the existing trade rows (3-6)
$ws.Cells.Item(7, 1).Value = 42649.64466435185   # A7 - Date (style 1 applied below)
$ws.Cells.Item(7, 2).Value = $false              # B7 - Profitable (boolean)
$ws.Cells.Item(7, 3).Value = 9923.7199999999993  # C7 - Principle
$ws.Cells.Item(7, 4).Value = 9980.61             # D7 - Start Principle
$ws.Cells.Item(7, 5).Value = 313.26998900000001  # E7 - BuyPrice
$ws.Cells.Item(7, 6).Value = 309.70001200000002  # F7 - SellPrice
$ws.Cells.Item(7, 7).Value = $false              # G7 - IsShortSell (boolean, style 1)
$ws.Cells.Item(7, 8).Value = -1.1399999999999999 # H7 - Price Change %
$ws.Cells.Item(7, 9).Value = $false              # I7 - Strong trade (boolean)

# Match the date-style formatting used by the other rows in columns A and G
# (use the canonical format code so it reuses the existing style instead of
# creating a brand new one)
$ws.Cells.Item(7, 1).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(7, 7).NumberFormat = "m/d/yy h:mm"

# Update column widths: columns E and F now share a single wider width
$ws.Columns.Item(5).ColumnWidth = 10
$ws.Columns.Item(6).ColumnWidth = 10
